$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.585.94"
$ws.Range("E2").Value = "  +6.96%  "
$ws.Range("D3").Value = "3.588.03"
$ws.Range("E3").Value = "  +2.86%  "
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'418.23"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "'130.50"
$ws.Range("E6").Value = "  -1.12%  "
$ws.Range("D7").Value = "'0.650"
$ws.Range("E7").Value = "  +2.98%  "
$ws.Range("D8").Value = "3.580.41"
$ws.Range("E8").Value = "  +2.85%  "
$ws.Range("D9").Value = "'0.998"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("D10").Value = "'0.769"
$ws.Range("E10").Value = "  +4.40%  "
$ws.Range("D11").Value = "'0.179"
$ws.Range("E11").Value = "  +16.03%  "
$ws.Range("D12").Value = "'0.0000344"
$ws.Range("E12").Value = "  +50.26%  "
$ws.Range("D13").Value = "'42.51"
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "4.135.32"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "'20.47"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "3.597.38"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("E19").Value = "  +4.92%  "
$ws.Range("D20").Value = "67.401.62"
$ws.Range("E20").Value = "  +6.84%  "
$ws.Range("D21").Value = "'12.38"
$ws.Range("E21").Value = "  -2.64%  "
$ws.Range("D22").Value = "'462.53"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").Value = "'88.54"
$ws.Range("E23").Value = "  -2.85%  "
$ws.Range("D24").Value = "'3.13"
$ws.Range("E24").Value = "  -5.58%  "
$ws.Range("D25").Value = "'13.43"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("D26").Value = "'3.37"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("D27").Value = "'10.24"
$ws.Range("E27").Value = "  -5.20%  "
$ws.Range("D28").Value = "'35.37"
$ws.Range("E28").Value = "  +4.81%  "
$ws.Range("D29").Value = "'4.84"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("E30").Value = "  +4.42%  "
$ws.Range("D31").Value = "'12.47"
$ws.Range("E31").Value = "  +2.19%  "
$ws.Range("E32").Value = "  +4.50%  "
$ws.Range("D33").Value = "'7.45"
$ws.Range("E33").Value = "  -2.02%  "
$ws.Range("D34").Value = "'41.70"
$ws.Range("E34").Value = "  +1.04%  "
$ws.Range("E35").Value = "  -4.15%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("D37").Value = "'56.82"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("D38").Value = "'0.0495"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").Value = "0.0₃0719"
$ws.Range("E39").Value = "  +21.80%  "
$ws.Range("E40").Value = "  +7.68%  "
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("D42").Value = "'3.04"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("D43").Value = "'148.62"
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "'3.27"
$ws.Range("E45").Value = "  -2.49%  "
$ws.Range("D46").Value = "'4.33"
$ws.Range("E46").Value = "  -3.67%  "
$ws.Range("D47").Value = "'0.312"
$ws.Range("E47").Value = "  -3.41%  "
$ws.Range("D48").Value = "'1.98"
$ws.Range("E48").Value = "  -4.57%  "
$ws.Range("D49").Value = "'2.34"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("D50").Value = "'2.72"
$ws.Range("E50").Value = "  +16.81%  "
$ws.Range("D51").Value = "'15.72"
$ws.Range("E51").Value = "  -4.74%  "
